# fix mk12 upper again
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("m4-uppers")

# Update the shared formula covering N3:N14 so every cell recalculates
# with J/200 instead of J/300. Setting the Formula on the whole range at
# once keeps Excel's shared-formula grouping intact on save.
$ws.Range("N3:N14").Formula = "=C3-D3*20-E3*0.8-F3*0.6-H3*5+J3/200"

# Update the standalone formula in N15 the same way.
$ws.Range("N15").Formula = "=C15-D15*20-E15*0.8-F15*0.6-H15*5+J15/200"

# Row 14 values: E14 2 -> 1, F14 0 -> 2
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 2

# Update the selected cell shown when the sheet is viewed.
$ws.Range("G10").Select()
